$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Is Available" column (E) was storing text "True"/"False" as strings.
# Push up set value: store real Boolean values instead.
$ws.Range("E2").Value = $true
$ws.Range("E3").Value = $false
